$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Cells.Item(2,6).Value = 32782.94
$ws1.Cells.Item(2,7).Value = 4
$ws1.Cells.Item(2,8).Value = 4
$ws2.Cells.Item(2,2).Value = 32782.94
$ws2.Cells.Item(2,3).Value = 1
$ws2.Cells.Item(2,4).Value = 2
$ws1.Cells.Item(3,6).Value = 64783.93
$ws1.Cells.Item(3,7).Value = 5
$ws1.Cells.Item(3,8).Value = 3
$ws2.Cells.Item(3,2).Value = 64783.93
$ws2.Cells.Item(3,3).Value = 5
$ws2.Cells.Item(3,4).Value = 4
$ws1.Cells.Item(4,6).Value = 68844.25
$ws1.Cells.Item(4,7).Value = 10
$ws1.Cells.Item(4,8).Value = 1
$ws2.Cells.Item(4,2).Value = 68844.25
$ws2.Cells.Item(4,3).Value = 4
$ws2.Cells.Item(4,4).Value = 5
$ws1.Cells.Item(5,6).Value = 92570.64
$ws1.Cells.Item(5,7).Value = 1
$ws1.Cells.Item(5,8).Value = 4
$ws2.Cells.Item(5,2).Value = 92570.64
$ws2.Cells.Item(5,3).Value = 10
$ws2.Cells.Item(5,4).Value = 5
$ws1.Cells.Item(6,6).Value = 76810.39999999999
$ws1.Cells.Item(6,7).Value = 5
$ws1.Cells.Item(6,8).Value = 3
$ws2.Cells.Item(6,2).Value = 76810.39999999999
$ws2.Cells.Item(6,3).Value = 8
$ws2.Cells.Item(6,4).Value = 2
$ws1.Cells.Item(7,6).Value = 77706.64
$ws1.Cells.Item(7,7).Value = 3
$ws1.Cells.Item(7,8).Value = 4
$ws2.Cells.Item(7,2).Value = 77706.64
$ws2.Cells.Item(7,3).Value = 10
$ws2.Cells.Item(7,4).Value = 1
$ws1.Cells.Item(8,6).Value = 93011.63
$ws1.Cells.Item(8,7).Value = 3
$ws1.Cells.Item(8,8).Value = 2
$ws2.Cells.Item(8,2).Value = 93011.63
$ws2.Cells.Item(8,3).Value = 6
$ws2.Cells.Item(8,4).Value = 2
$ws1.Cells.Item(9,6).Value = 87330.58
$ws1.Cells.Item(9,7).Value = 9
$ws1.Cells.Item(9,8).Value = 4
$ws2.Cells.Item(9,2).Value = 87330.58
$ws2.Cells.Item(9,3).Value = 6
$ws2.Cells.Item(9,4).Value = 2
$ws1.Cells.Item(10,6).Value = 51216.11
$ws1.Cells.Item(10,7).Value = 7
$ws1.Cells.Item(10,8).Value = 5
$ws2.Cells.Item(10,2).Value = 51216.11
$ws2.Cells.Item(10,3).Value = 4
$ws2.Cells.Item(10,4).Value = 2
$ws1.Cells.Item(11,6).Value = 54035.1
$ws1.Cells.Item(11,7).Value = 9
$ws1.Cells.Item(11,8).Value = 1
$ws2.Cells.Item(11,2).Value = 54035.1
$ws2.Cells.Item(11,3).Value = 9
$ws2.Cells.Item(11,4).Value = 1
$ws1.Cells.Item(12,6).Value = 84947.84
$ws1.Cells.Item(12,7).Value = 6
$ws1.Cells.Item(12,8).Value = 1
$ws2.Cells.Item(12,2).Value = 84947.84
$ws2.Cells.Item(12,3).Value = 8
$ws2.Cells.Item(12,4).Value = 5
$ws1.Cells.Item(13,6).Value = 46564.92
$ws1.Cells.Item(13,7).Value = 3
$ws1.Cells.Item(13,8).Value = 5
$ws2.Cells.Item(13,2).Value = 46564.92
$ws2.Cells.Item(13,3).Value = 7
$ws2.Cells.Item(13,4).Value = 5
$ws1.Cells.Item(14,6).Value = 36589.02
$ws1.Cells.Item(14,7).Value = 5
$ws1.Cells.Item(14,8).Value = 5
$ws2.Cells.Item(14,2).Value = 36589.02
$ws2.Cells.Item(14,3).Value = 4
$ws2.Cells.Item(14,4).Value = 3
$ws1.Cells.Item(15,6).Value = 52951.43
$ws1.Cells.Item(15,7).Value = 2
$ws1.Cells.Item(15,8).Value = 3
$ws2.Cells.Item(15,2).Value = 52951.43
$ws2.Cells.Item(15,3).Value = 5
$ws2.Cells.Item(15,4).Value = 3
$ws1.Cells.Item(16,6).Value = 34252.34
$ws1.Cells.Item(16,7).Value = 7
$ws1.Cells.Item(16,8).Value = 3
$ws2.Cells.Item(16,2).Value = 34252.34
$ws2.Cells.Item(16,3).Value = 2
$ws2.Cells.Item(16,4).Value = 2
$ws1.Cells.Item(17,6).Value = 35251.46
$ws1.Cells.Item(17,7).Value = 5
$ws1.Cells.Item(17,8).Value = 5
$ws2.Cells.Item(17,2).Value = 35251.46
$ws2.Cells.Item(17,3).Value = 1
$ws2.Cells.Item(17,4).Value = 4
$ws1.Cells.Item(18,6).Value = 53432.94
$ws1.Cells.Item(18,7).Value = 6
$ws1.Cells.Item(18,8).Value = 5
$ws2.Cells.Item(18,2).Value = 53432.94
$ws2.Cells.Item(18,3).Value = 9
$ws2.Cells.Item(18,4).Value = 3
$ws1.Cells.Item(19,6).Value = 93101.91
$ws1.Cells.Item(19,7).Value = 8
$ws1.Cells.Item(19,8).Value = 2
$ws2.Cells.Item(19,2).Value = 93101.91
$ws2.Cells.Item(19,3).Value = 1
$ws2.Cells.Item(19,4).Value = 3
$ws1.Cells.Item(20,6).Value = 46310.09
$ws1.Cells.Item(20,7).Value = 7
$ws1.Cells.Item(20,8).Value = 2
$ws2.Cells.Item(20,2).Value = 46310.09
$ws2.Cells.Item(20,3).Value = 9
$ws2.Cells.Item(20,4).Value = 3
$ws1.Cells.Item(21,6).Value = 69783.53
$ws1.Cells.Item(21,7).Value = 3
$ws1.Cells.Item(21,8).Value = 3
$ws2.Cells.Item(21,2).Value = 69783.53
$ws2.Cells.Item(21,3).Value = 2
$ws2.Cells.Item(21,4).Value = 5
$ws1.Cells.Item(22,6).Value = 62239.46
$ws1.Cells.Item(22,7).Value = 4
$ws1.Cells.Item(22,8).Value = 2
$ws2.Cells.Item(22,2).Value = 62239.46
$ws2.Cells.Item(22,3).Value = 2
$ws2.Cells.Item(22,4).Value = 5
$ws1.Cells.Item(23,6).Value = 66203.67
$ws1.Cells.Item(23,7).Value = 2
$ws1.Cells.Item(23,8).Value = 3
$ws2.Cells.Item(23,2).Value = 66203.67
$ws2.Cells.Item(23,3).Value = 8
$ws2.Cells.Item(23,4).Value = 3
$ws1.Cells.Item(24,6).Value = 33779.86
$ws1.Cells.Item(24,7).Value = 6
$ws1.Cells.Item(24,8).Value = 4
$ws2.Cells.Item(24,2).Value = 33779.86
$ws2.Cells.Item(24,3).Value = 3
$ws2.Cells.Item(24,4).Value = 1
$ws1.Cells.Item(25,6).Value = 68757.35000000001
$ws1.Cells.Item(25,7).Value = 1
$ws1.Cells.Item(25,8).Value = 1
$ws2.Cells.Item(25,2).Value = 68757.35000000001
$ws2.Cells.Item(25,3).Value = 9
$ws2.Cells.Item(25,4).Value = 5
$ws1.Cells.Item(26,6).Value = 68970.78999999999
$ws1.Cells.Item(26,7).Value = 7
$ws1.Cells.Item(26,8).Value = 4
$ws2.Cells.Item(26,2).Value = 68970.78999999999
$ws2.Cells.Item(26,3).Value = 1
$ws2.Cells.Item(26,4).Value = 4
$ws1.Cells.Item(27,6).Value = 66939.73
$ws1.Cells.Item(27,7).Value = 7
$ws1.Cells.Item(27,8).Value = 4
$ws2.Cells.Item(27,2).Value = 66939.73
$ws2.Cells.Item(27,3).Value = 4
$ws2.Cells.Item(27,4).Value = 1
$ws1.Cells.Item(28,6).Value = 82147.63
$ws1.Cells.Item(28,7).Value = 1
$ws1.Cells.Item(28,8).Value = 5
$ws2.Cells.Item(28,2).Value = 82147.63
$ws2.Cells.Item(28,3).Value = 5
$ws2.Cells.Item(28,4).Value = 2
$ws1.Cells.Item(29,6).Value = 68103.28999999999
$ws1.Cells.Item(29,7).Value = 4
$ws1.Cells.Item(29,8).Value = 5
$ws2.Cells.Item(29,2).Value = 68103.28999999999
$ws2.Cells.Item(29,3).Value = 3
$ws2.Cells.Item(29,4).Value = 4
$ws1.Cells.Item(30,6).Value = 66163.17999999999
$ws1.Cells.Item(30,7).Value = 3
$ws1.Cells.Item(30,8).Value = 1
$ws2.Cells.Item(30,2).Value = 66163.17999999999
$ws2.Cells.Item(30,3).Value = 5
$ws2.Cells.Item(30,4).Value = 4
$ws1.Cells.Item(31,6).Value = 80505.92
$ws1.Cells.Item(31,7).Value = 2
$ws1.Cells.Item(31,8).Value = 3
$ws2.Cells.Item(31,2).Value = 80505.92
$ws2.Cells.Item(31,3).Value = 3
$ws2.Cells.Item(31,4).Value = 5
$ws1.Cells.Item(32,6).Value = 92134.13
$ws1.Cells.Item(32,7).Value = 10
$ws1.Cells.Item(32,8).Value = 4
$ws2.Cells.Item(32,2).Value = 92134.13
$ws2.Cells.Item(32,3).Value = 4
$ws2.Cells.Item(32,4).Value = 2
$ws1.Cells.Item(33,6).Value = 77663.71000000001
$ws1.Cells.Item(33,7).Value = 3
$ws1.Cells.Item(33,8).Value = 5
$ws2.Cells.Item(33,2).Value = 77663.71000000001
$ws2.Cells.Item(33,3).Value = 8
$ws2.Cells.Item(33,4).Value = 5
$ws1.Cells.Item(34,6).Value = 93058.49000000001
$ws1.Cells.Item(34,7).Value = 7
$ws1.Cells.Item(34,8).Value = 3
$ws2.Cells.Item(34,2).Value = 93058.49000000001
$ws2.Cells.Item(34,3).Value = 6
$ws2.Cells.Item(34,4).Value = 1
$ws1.Cells.Item(35,6).Value = 48891.26
$ws1.Cells.Item(35,7).Value = 6
$ws1.Cells.Item(35,8).Value = 4
$ws2.Cells.Item(35,2).Value = 48891.26
$ws2.Cells.Item(35,3).Value = 7
$ws2.Cells.Item(35,4).Value = 2
$ws1.Cells.Item(36,6).Value = 64729.06
$ws1.Cells.Item(36,7).Value = 9
$ws1.Cells.Item(36,8).Value = 3
$ws2.Cells.Item(36,2).Value = 64729.06
$ws2.Cells.Item(36,3).Value = 10
$ws2.Cells.Item(36,4).Value = 1
$ws1.Cells.Item(37,6).Value = 40100.87
$ws1.Cells.Item(37,7).Value = 4
$ws1.Cells.Item(37,8).Value = 3
$ws2.Cells.Item(37,2).Value = 40100.87
$ws2.Cells.Item(37,3).Value = 4
$ws2.Cells.Item(37,4).Value = 5
$ws1.Cells.Item(38,6).Value = 65327.75
$ws1.Cells.Item(38,7).Value = 10
$ws1.Cells.Item(38,8).Value = 2
$ws2.Cells.Item(38,2).Value = 65327.75
$ws2.Cells.Item(38,3).Value = 4
$ws2.Cells.Item(38,4).Value = 3
$ws1.Cells.Item(39,6).Value = 70823.73
$ws1.Cells.Item(39,7).Value = 3
$ws1.Cells.Item(39,8).Value = 1
$ws2.Cells.Item(39,2).Value = 70823.73
$ws2.Cells.Item(39,3).Value = 3
$ws2.Cells.Item(39,4).Value = 1
$ws1.Cells.Item(40,6).Value = 72342.36
$ws1.Cells.Item(40,7).Value = 4
$ws1.Cells.Item(40,8).Value = 2
$ws2.Cells.Item(40,2).Value = 72342.36
$ws2.Cells.Item(40,3).Value = 1
$ws2.Cells.Item(40,4).Value = 4
$ws1.Cells.Item(41,6).Value = 73821.23
$ws1.Cells.Item(41,7).Value = 8
$ws1.Cells.Item(41,8).Value = 3
$ws2.Cells.Item(41,2).Value = 73821.23
$ws2.Cells.Item(41,3).Value = 2
$ws2.Cells.Item(41,4).Value = 5
$ws1.Cells.Item(42,6).Value = 43829.72
$ws1.Cells.Item(42,7).Value = 6
$ws1.Cells.Item(42,8).Value = 2
$ws2.Cells.Item(42,2).Value = 43829.72
$ws2.Cells.Item(42,3).Value = 10
$ws2.Cells.Item(42,4).Value = 3
$ws1.Cells.Item(43,6).Value = 73465.78
$ws1.Cells.Item(43,7).Value = 2
$ws1.Cells.Item(43,8).Value = 3
$ws2.Cells.Item(43,2).Value = 73465.78
$ws2.Cells.Item(43,3).Value = 3
$ws2.Cells.Item(43,4).Value = 3
$ws1.Cells.Item(44,6).Value = 39704.67
$ws1.Cells.Item(44,7).Value = 2
$ws1.Cells.Item(44,8).Value = 5
$ws2.Cells.Item(44,2).Value = 39704.67
$ws2.Cells.Item(44,3).Value = 1
$ws2.Cells.Item(44,4).Value = 4
$ws1.Cells.Item(45,6).Value = 70345.73
$ws1.Cells.Item(45,7).Value = 5
$ws1.Cells.Item(45,8).Value = 5
$ws2.Cells.Item(45,2).Value = 70345.73
$ws2.Cells.Item(45,3).Value = 4
$ws2.Cells.Item(45,4).Value = 3
$ws1.Cells.Item(46,6).Value = 81879.17
$ws1.Cells.Item(46,7).Value = 8
$ws1.Cells.Item(46,8).Value = 5
$ws2.Cells.Item(46,2).Value = 81879.17
$ws2.Cells.Item(46,3).Value = 9
$ws2.Cells.Item(46,4).Value = 4
$ws1.Cells.Item(47,6).Value = 97245.82000000001
$ws1.Cells.Item(47,7).Value = 10
$ws1.Cells.Item(47,8).Value = 5
$ws2.Cells.Item(47,2).Value = 97245.82000000001
$ws2.Cells.Item(47,3).Value = 3
$ws2.Cells.Item(47,4).Value = 1
$ws1.Cells.Item(48,6).Value = 79261.56
$ws1.Cells.Item(48,7).Value = 6
$ws1.Cells.Item(48,8).Value = 4
$ws2.Cells.Item(48,2).Value = 79261.56
$ws2.Cells.Item(48,3).Value = 4
$ws2.Cells.Item(48,4).Value = 3
$ws1.Cells.Item(49,6).Value = 47033.56
$ws1.Cells.Item(49,7).Value = 6
$ws1.Cells.Item(49,8).Value = 4
$ws2.Cells.Item(49,2).Value = 47033.56
$ws2.Cells.Item(49,3).Value = 8
$ws2.Cells.Item(49,4).Value = 4
$ws1.Cells.Item(50,6).Value = 73037.24000000001
$ws1.Cells.Item(50,7).Value = 6
$ws1.Cells.Item(50,8).Value = 3
$ws2.Cells.Item(50,2).Value = 73037.24000000001
$ws2.Cells.Item(50,3).Value = 1
$ws2.Cells.Item(50,4).Value = 2
$ws1.Cells.Item(51,6).Value = 38587.91
$ws1.Cells.Item(51,7).Value = 5
$ws1.Cells.Item(51,8).Value = 5
$ws2.Cells.Item(51,2).Value = 38587.91
$ws2.Cells.Item(51,3).Value = 4
$ws2.Cells.Item(51,4).Value = 5
$ws1.Cells.Item(52,6).Value = 65352.76
$ws1.Cells.Item(52,7).Value = 6
$ws1.Cells.Item(52,8).Value = 4
$ws2.Cells.Item(52,2).Value = 65352.76
$ws2.Cells.Item(52,3).Value = 2
$ws2.Cells.Item(52,4).Value = 4
$ws1.Cells.Item(53,6).Value = 44261.35
$ws1.Cells.Item(53,7).Value = 2
$ws1.Cells.Item(53,8).Value = 3
$ws2.Cells.Item(53,2).Value = 44261.35
$ws2.Cells.Item(53,3).Value = 10
$ws2.Cells.Item(53,4).Value = 5
$ws1.Cells.Item(54,6).Value = 36295.76
$ws1.Cells.Item(54,7).Value = 4
$ws1.Cells.Item(54,8).Value = 2
$ws2.Cells.Item(54,2).Value = 36295.76
$ws2.Cells.Item(54,3).Value = 4
$ws2.Cells.Item(54,4).Value = 5
$ws1.Cells.Item(55,6).Value = 67410.72
$ws1.Cells.Item(55,7).Value = 7
$ws1.Cells.Item(55,8).Value = 1
$ws2.Cells.Item(55,2).Value = 67410.72
$ws2.Cells.Item(55,3).Value = 7
$ws2.Cells.Item(55,4).Value = 4
$ws1.Cells.Item(56,6).Value = 30463.77
$ws1.Cells.Item(56,7).Value = 9
$ws1.Cells.Item(56,8).Value = 5
$ws2.Cells.Item(56,2).Value = 30463.77
$ws2.Cells.Item(56,3).Value = 8
$ws2.Cells.Item(56,4).Value = 3
$ws1.Cells.Item(57,6).Value = 37704.98
$ws1.Cells.Item(57,7).Value = 7
$ws1.Cells.Item(57,8).Value = 4
$ws2.Cells.Item(57,2).Value = 37704.98
$ws2.Cells.Item(57,3).Value = 8
$ws2.Cells.Item(57,4).Value = 1
$ws1.Cells.Item(58,6).Value = 45808.88
$ws1.Cells.Item(58,7).Value = 10
$ws1.Cells.Item(58,8).Value = 3
$ws2.Cells.Item(58,2).Value = 45808.88
$ws2.Cells.Item(58,3).Value = 10
$ws2.Cells.Item(58,4).Value = 3
$ws1.Cells.Item(59,6).Value = 64630.3
$ws1.Cells.Item(59,7).Value = 7
$ws1.Cells.Item(59,8).Value = 5
$ws2.Cells.Item(59,2).Value = 64630.3
$ws2.Cells.Item(59,3).Value = 6
$ws2.Cells.Item(59,4).Value = 2
$ws1.Cells.Item(60,6).Value = 88981.64999999999
$ws1.Cells.Item(60,7).Value = 10
$ws1.Cells.Item(60,8).Value = 2
$ws2.Cells.Item(60,2).Value = 88981.64999999999
$ws2.Cells.Item(60,3).Value = 9
$ws2.Cells.Item(60,4).Value = 1
$ws1.Cells.Item(61,6).Value = 86762.95
$ws1.Cells.Item(61,7).Value = 5
$ws1.Cells.Item(61,8).Value = 4
$ws2.Cells.Item(61,2).Value = 86762.95
$ws2.Cells.Item(61,3).Value = 3
$ws2.Cells.Item(61,4).Value = 5
$ws1.Cells.Item(62,6).Value = 55374.25
$ws1.Cells.Item(62,7).Value = 3
$ws1.Cells.Item(62,8).Value = 4
$ws2.Cells.Item(62,2).Value = 55374.25
$ws2.Cells.Item(62,3).Value = 4
$ws2.Cells.Item(62,4).Value = 4
$ws1.Cells.Item(63,6).Value = 44478.14
$ws1.Cells.Item(63,7).Value = 3
$ws1.Cells.Item(63,8).Value = 2
$ws2.Cells.Item(63,2).Value = 44478.14
$ws2.Cells.Item(63,3).Value = 6
$ws2.Cells.Item(63,4).Value = 5
$ws1.Cells.Item(64,6).Value = 44659.69
$ws1.Cells.Item(64,7).Value = 4
$ws1.Cells.Item(64,8).Value = 5
$ws2.Cells.Item(64,2).Value = 44659.69
$ws2.Cells.Item(64,3).Value = 1
$ws2.Cells.Item(64,4).Value = 4
$ws1.Cells.Item(65,6).Value = 50189.87
$ws1.Cells.Item(65,7).Value = 9
$ws1.Cells.Item(65,8).Value = 1
$ws2.Cells.Item(65,2).Value = 50189.87
$ws2.Cells.Item(65,3).Value = 4
$ws2.Cells.Item(65,4).Value = 5
$ws1.Cells.Item(66,6).Value = 85888.91
$ws1.Cells.Item(66,7).Value = 4
$ws1.Cells.Item(66,8).Value = 5
$ws2.Cells.Item(66,2).Value = 85888.91
$ws2.Cells.Item(66,3).Value = 8
$ws2.Cells.Item(66,4).Value = 5
$ws1.Cells.Item(67,6).Value = 71546.86
$ws1.Cells.Item(67,7).Value = 7
$ws1.Cells.Item(67,8).Value = 1
$ws2.Cells.Item(67,2).Value = 71546.86
$ws2.Cells.Item(67,3).Value = 6
$ws2.Cells.Item(67,4).Value = 1
$ws1.Cells.Item(68,6).Value = 49673.93
$ws1.Cells.Item(68,7).Value = 3
$ws1.Cells.Item(68,8).Value = 2
$ws2.Cells.Item(68,2).Value = 49673.93
$ws2.Cells.Item(68,3).Value = 4
$ws2.Cells.Item(68,4).Value = 2
$ws1.Cells.Item(69,6).Value = 52075.79
$ws1.Cells.Item(69,7).Value = 9
$ws1.Cells.Item(69,8).Value = 5
$ws2.Cells.Item(69,2).Value = 52075.79
$ws2.Cells.Item(69,3).Value = 8
$ws2.Cells.Item(69,4).Value = 5
$ws1.Cells.Item(70,6).Value = 95359.23
$ws1.Cells.Item(70,7).Value = 9
$ws1.Cells.Item(70,8).Value = 3
$ws2.Cells.Item(70,2).Value = 95359.23
$ws2.Cells.Item(70,3).Value = 10
$ws2.Cells.Item(70,4).Value = 1
$ws1.Cells.Item(71,6).Value = 89207.82000000001
$ws1.Cells.Item(71,7).Value = 4
$ws1.Cells.Item(71,8).Value = 3
$ws2.Cells.Item(71,2).Value = 89207.82000000001
$ws2.Cells.Item(71,3).Value = 1
$ws2.Cells.Item(71,4).Value = 1
$ws1.Cells.Item(72,6).Value = 86615.09
$ws1.Cells.Item(72,7).Value = 5
$ws1.Cells.Item(72,8).Value = 2
$ws2.Cells.Item(72,2).Value = 86615.09
$ws2.Cells.Item(72,3).Value = 10
$ws2.Cells.Item(72,4).Value = 1
$ws1.Cells.Item(73,6).Value = 86711.21000000001
$ws1.Cells.Item(73,7).Value = 2
$ws1.Cells.Item(73,8).Value = 5
$ws2.Cells.Item(73,2).Value = 86711.21000000001
$ws2.Cells.Item(73,3).Value = 8
$ws2.Cells.Item(73,4).Value = 4
$ws1.Cells.Item(74,6).Value = 58254.1
$ws1.Cells.Item(74,7).Value = 1
$ws1.Cells.Item(74,8).Value = 5
$ws2.Cells.Item(74,2).Value = 58254.1
$ws2.Cells.Item(74,3).Value = 3
$ws2.Cells.Item(74,4).Value = 3
$ws1.Cells.Item(75,6).Value = 59893.33
$ws1.Cells.Item(75,7).Value = 10
$ws1.Cells.Item(75,8).Value = 5
$ws2.Cells.Item(75,2).Value = 59893.33
$ws2.Cells.Item(75,3).Value = 5
$ws2.Cells.Item(75,4).Value = 3
$ws1.Cells.Item(76,6).Value = 30389.22
$ws1.Cells.Item(76,7).Value = 8
$ws1.Cells.Item(76,8).Value = 5
$ws2.Cells.Item(76,2).Value = 30389.22
$ws2.Cells.Item(76,3).Value = 8
$ws2.Cells.Item(76,4).Value = 1
$ws1.Cells.Item(77,6).Value = 46985.51
$ws1.Cells.Item(77,7).Value = 9
$ws1.Cells.Item(77,8).Value = 3
$ws2.Cells.Item(77,2).Value = 46985.51
$ws2.Cells.Item(77,3).Value = 5
$ws2.Cells.Item(77,4).Value = 3
$ws1.Cells.Item(78,6).Value = 90303.08
$ws1.Cells.Item(78,7).Value = 3
$ws1.Cells.Item(78,8).Value = 3
$ws2.Cells.Item(78,2).Value = 90303.08
$ws2.Cells.Item(78,3).Value = 3
$ws2.Cells.Item(78,4).Value = 3
$ws1.Cells.Item(79,6).Value = 83861.98
$ws1.Cells.Item(79,7).Value = 4
$ws1.Cells.Item(79,8).Value = 5
$ws2.Cells.Item(79,2).Value = 83861.98
$ws2.Cells.Item(79,3).Value = 3
$ws2.Cells.Item(79,4).Value = 4
$ws1.Cells.Item(80,6).Value = 95595
$ws1.Cells.Item(80,7).Value = 8
$ws1.Cells.Item(80,8).Value = 5
$ws2.Cells.Item(80,2).Value = 95595
$ws2.Cells.Item(80,3).Value = 6
$ws2.Cells.Item(80,4).Value = 1
$ws1.Cells.Item(81,6).Value = 87752.85000000001
$ws1.Cells.Item(81,7).Value = 10
$ws1.Cells.Item(81,8).Value = 1
$ws2.Cells.Item(81,2).Value = 87752.85000000001
$ws2.Cells.Item(81,3).Value = 5
$ws2.Cells.Item(81,4).Value = 4
$ws1.Cells.Item(82,6).Value = 61722.04
$ws1.Cells.Item(82,7).Value = 5
$ws1.Cells.Item(82,8).Value = 2
$ws2.Cells.Item(82,2).Value = 61722.04
$ws2.Cells.Item(82,3).Value = 5
$ws2.Cells.Item(82,4).Value = 3
$ws1.Cells.Item(83,6).Value = 49537.94
$ws1.Cells.Item(83,7).Value = 5
$ws1.Cells.Item(83,8).Value = 1
$ws2.Cells.Item(83,2).Value = 49537.94
$ws2.Cells.Item(83,3).Value = 7
$ws2.Cells.Item(83,4).Value = 5
$ws1.Cells.Item(84,6).Value = 76340.49000000001
$ws1.Cells.Item(84,7).Value = 10
$ws1.Cells.Item(84,8).Value = 2
$ws2.Cells.Item(84,2).Value = 76340.49000000001
$ws2.Cells.Item(84,3).Value = 8
$ws2.Cells.Item(84,4).Value = 1
$ws1.Cells.Item(85,6).Value = 68742.64999999999
$ws1.Cells.Item(85,7).Value = 1
$ws1.Cells.Item(85,8).Value = 4
$ws2.Cells.Item(85,2).Value = 68742.64999999999
$ws2.Cells.Item(85,3).Value = 3
$ws2.Cells.Item(85,4).Value = 1
$ws1.Cells.Item(86,6).Value = 36911.73
$ws1.Cells.Item(86,7).Value = 6
$ws1.Cells.Item(86,8).Value = 4
$ws2.Cells.Item(86,2).Value = 36911.73
$ws2.Cells.Item(86,3).Value = 9
$ws2.Cells.Item(86,4).Value = 5
$ws1.Cells.Item(87,6).Value = 81207.81
$ws1.Cells.Item(87,7).Value = 4
$ws1.Cells.Item(87,8).Value = 4
$ws2.Cells.Item(87,2).Value = 81207.81
$ws2.Cells.Item(87,3).Value = 4
$ws2.Cells.Item(87,4).Value = 3
$ws1.Cells.Item(88,6).Value = 43349.55
$ws1.Cells.Item(88,7).Value = 10
$ws1.Cells.Item(88,8).Value = 5
$ws2.Cells.Item(88,2).Value = 43349.55
$ws2.Cells.Item(88,3).Value = 3
$ws2.Cells.Item(88,4).Value = 2
$ws1.Cells.Item(89,6).Value = 90515.27
$ws1.Cells.Item(89,7).Value = 8
$ws1.Cells.Item(89,8).Value = 4
$ws2.Cells.Item(89,2).Value = 90515.27
$ws2.Cells.Item(89,3).Value = 1
$ws2.Cells.Item(89,4).Value = 3
$ws1.Cells.Item(90,6).Value = 36473.26
$ws1.Cells.Item(90,7).Value = 9
$ws1.Cells.Item(90,8).Value = 1
$ws2.Cells.Item(90,2).Value = 36473.26
$ws2.Cells.Item(90,3).Value = 7
$ws2.Cells.Item(90,4).Value = 1
$ws1.Cells.Item(91,6).Value = 57243.91
$ws1.Cells.Item(91,7).Value = 5
$ws1.Cells.Item(91,8).Value = 2
$ws2.Cells.Item(91,2).Value = 57243.91
$ws2.Cells.Item(91,3).Value = 1
$ws2.Cells.Item(91,4).Value = 1
$ws1.Cells.Item(92,6).Value = 40745.12
$ws1.Cells.Item(92,7).Value = 6
$ws1.Cells.Item(92,8).Value = 2
$ws2.Cells.Item(92,2).Value = 40745.12
$ws2.Cells.Item(92,3).Value = 2
$ws2.Cells.Item(92,4).Value = 2
$ws1.Cells.Item(93,6).Value = 54512.92
$ws1.Cells.Item(93,7).Value = 1
$ws1.Cells.Item(93,8).Value = 1
$ws2.Cells.Item(93,2).Value = 54512.92
$ws2.Cells.Item(93,3).Value = 10
$ws2.Cells.Item(93,4).Value = 2
$ws1.Cells.Item(94,6).Value = 35942.48
$ws1.Cells.Item(94,7).Value = 9
$ws1.Cells.Item(94,8).Value = 4
$ws2.Cells.Item(94,2).Value = 35942.48
$ws2.Cells.Item(94,3).Value = 6
$ws2.Cells.Item(94,4).Value = 2
$ws1.Cells.Item(95,6).Value = 51157.87
$ws1.Cells.Item(95,7).Value = 1
$ws1.Cells.Item(95,8).Value = 2
$ws2.Cells.Item(95,2).Value = 51157.87
$ws2.Cells.Item(95,3).Value = 7
$ws2.Cells.Item(95,4).Value = 5
$ws1.Cells.Item(96,6).Value = 46048.21
$ws1.Cells.Item(96,7).Value = 5
$ws1.Cells.Item(96,8).Value = 2
$ws2.Cells.Item(96,2).Value = 46048.21
$ws2.Cells.Item(96,3).Value = 2
$ws2.Cells.Item(96,4).Value = 4
$ws1.Cells.Item(97,6).Value = 94677.07000000001
$ws1.Cells.Item(97,7).Value = 8
$ws1.Cells.Item(97,8).Value = 4
$ws2.Cells.Item(97,2).Value = 94677.07000000001
$ws2.Cells.Item(97,3).Value = 1
$ws2.Cells.Item(97,4).Value = 4
$ws1.Cells.Item(98,6).Value = 62024.31
$ws1.Cells.Item(98,7).Value = 6
$ws1.Cells.Item(98,8).Value = 4
$ws2.Cells.Item(98,2).Value = 62024.31
$ws2.Cells.Item(98,3).Value = 7
$ws2.Cells.Item(98,4).Value = 4
$ws1.Cells.Item(99,6).Value = 68702.19
$ws1.Cells.Item(99,7).Value = 3
$ws1.Cells.Item(99,8).Value = 2
$ws2.Cells.Item(99,2).Value = 68702.19
$ws2.Cells.Item(99,3).Value = 9
$ws2.Cells.Item(99,4).Value = 4
$ws1.Cells.Item(100,6).Value = 74990.98
$ws1.Cells.Item(100,7).Value = 6
$ws1.Cells.Item(100,8).Value = 5
$ws2.Cells.Item(100,2).Value = 74990.98
$ws2.Cells.Item(100,3).Value = 9
$ws2.Cells.Item(100,4).Value = 4
$ws1.Cells.Item(101,6).Value = 62288.6
$ws1.Cells.Item(101,7).Value = 3
$ws1.Cells.Item(101,8).Value = 2
$ws2.Cells.Item(101,2).Value = 62288.6
$ws2.Cells.Item(101,3).Value = 2
$ws2.Cells.Item(101,4).Value = 3
